$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C9").Value = "test"
$ws.Range("C10").Select()
